$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 145, shifting existing rows 145:246 down to 146:247.
$ws.Rows(145).Insert()

# Populate the newly inserted row 145 with the new data record.
$ws.Range("A145").Value = 11
$ws.Range("B145").Value = "Vega Monumental Concepción"
$ws.Range("C145").Value = "Bíobío"
$ws.Range("D145").Value = 44651
$ws.Range("E145").Value = 8
$ws.Range("F145").Value = "Fruta"
$ws.Range("G145").Value = 100102
$ws.Range("H145").Value = "Cítricos"
$ws.Range("I145").Value = 100102005
$ws.Range("J145").Value = "Naranja"
$ws.Range("K145").Value = "Valencia"
$ws.Range("L145").Value = "Primera"
$ws.Range("M145").Value = 270
$ws.Range("N145").Value = 8500
$ws.Range("O145").Value = 9000
$ws.Range("P145").Value = 8722
$ws.Range("Q145").Value = "$/caja 15 kilos granel"
$ws.Range("R145").Value = "Región de O'Higgins"
$ws.Range("S145").Value = 581
$ws.Range("T145").Value = 15
